# Weekly update: new crime data collected (report period shifted one week forward)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: Volume/Number and reporting week range ---
# "Volume 32   Number  24" -> "...25"  (rich-text run; substring replace keeps formatting)
$ws.Range("A8").Characters(21, 2).Text = "25"
# "Report Covering the Week  6/9/2025  Through  6/15/2025" -> 6/16/2025 .. 6/22/2025
# Replace the later (rightmost) date first so the earlier date's position is unaffected
# by the one-character length change of the first date (8 -> 9 chars).
$ws.Range("C9").Characters(46, 9).Text = "6/22/2025"
$ws.Range("C9").Characters(27, 8).Text = "6/16/2025"

# --- Crime statistics table (rows 14-30), columns C:N: only cells that actually change ---
# Row 14
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 0
$ws.Range("N14").Value = -84.615384615384

# Row 15
$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 40
$ws.Range("I15").Value = 22
$ws.Range("K15").Value = 57.142857142857
$ws.Range("L15").Value = 69.230769230769
$ws.Range("M15").Value = 144.444444444444
$ws.Range("N15").Value = -35.294117647058

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 49
$ws.Range("H16").Value = -42.857142857142
$ws.Range("I16").Value = 154
$ws.Range("J16").Value = 195
$ws.Range("K16").Value = -21.025641025641
$ws.Range("L16").Value = -4.938271604938
$ws.Range("M16").Value = -31.25
$ws.Range("N16").Value = -79.575596816976

# Row 17
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = -15
$ws.Range("F17").Value = 64
$ws.Range("G17").Value = 59
$ws.Range("H17").Value = 8.474576271186
$ws.Range("I17").Value = 356
$ws.Range("J17").Value = 346
$ws.Range("K17").Value = 2.890173410404
$ws.Range("L17").Value = 10.559006211180
$ws.Range("M17").Value = 145.51724137931
$ws.Range("N17").Value = -8.717948717948

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -37.5
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 107
$ws.Range("J18").Value = 87
$ws.Range("K18").Value = 22.988505747126
$ws.Range("L18").Value = 37.179487179487
$ws.Range("M18").Value = -13.709677419354
$ws.Range("N18").Value = -80.330882352941

# Row 19
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 13
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -10.638297872340
$ws.Range("I19").Value = 243
$ws.Range("J19").Value = 298
$ws.Range("K19").Value = -18.456375838926
$ws.Range("L19").Value = -5.813953488372
$ws.Range("M19").Value = 36.516853932584
$ws.Range("N19").Value = -55.248618784530

# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 133.333333333333
$ws.Range("F20").Value = 17
$ws.Range("H20").Value = 30.769230769230
$ws.Range("I20").Value = 95
$ws.Range("J20").Value = 95
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -15.178571428571
$ws.Range("M20").Value = 10.465116279069
$ws.Range("N20").Value = -85.925925925925

# Row 21
$ws.Range("C21").Value = 50
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = -9.090909090909
$ws.Range("F21").Value = 178
$ws.Range("G21").Value = 189
$ws.Range("H21").Value = -5.820105820105
$ws.Range("I21").Value = 979
$ws.Range("J21").Value = 1037
$ws.Range("K21").Value = -5.593056894889
$ws.Range("L21").Value = 3.597883597883
$ws.Range("M21").Value = 26.977950713359
$ws.Range("N21").Value = -66.847273958686

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("I22").Value = 18
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = 12.5
$ws.Range("L22").Value = 80
$ws.Range("M22").Value = 63.636363636363

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = "'0"
$ws.Range("E23").Value = "***.*"
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 30
$ws.Range("K23").Value = 20
$ws.Range("L23").Value = 66.666666666666
$ws.Range("M23").Value = 57.894736842105

# Row 24
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 62
$ws.Range("E24").Value = -62.903225806451
$ws.Range("F24").Value = 135
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = -19.161676646706
$ws.Range("I24").Value = 913
$ws.Range("J24").Value = 866
$ws.Range("K24").Value = 5.427251732101
$ws.Range("L24").Value = 22.550335570469
$ws.Range("M24").Value = 65.698729582577

# Row 25
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 41
$ws.Range("E25").Value = -68.292682926829
$ws.Range("F25").Value = 69
$ws.Range("G25").Value = 123
$ws.Range("H25").Value = -43.902439024390
$ws.Range("I25").Value = 514
$ws.Range("J25").Value = 543
$ws.Range("K25").Value = -5.340699815837
$ws.Range("L25").Value = 65.273311897106

# Row 26
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = -36
$ws.Range("F26").Value = 71
$ws.Range("G26").Value = 101
$ws.Range("H26").Value = -29.702970297029
$ws.Range("I26").Value = 419
$ws.Range("J26").Value = 482
$ws.Range("K26").Value = -13.070539419087
$ws.Range("L26").Value = -9.892473118279
$ws.Range("M26").Value = 10.263157894736

# Row 27
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 23
$ws.Range("K27").Value = 4.545454545454
$ws.Range("L27").Value = 9.523809523809

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -7.692307692307
$ws.Range("I28").Value = 57
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = -5
$ws.Range("L28").Value = 42.5

# Row 29
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "***.*"
$ws.Range("M29").Value = -45.454545454545
$ws.Range("N29").Value = -83.333333333333

# Row 30
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "***.*"
$ws.Range("M30").Value = -76.923076923076
$ws.Range("N30").Value = -95.588235294117
